# Scheduled-runner style update of market-price derived columns (H:N)
# across the per-job Leve sheets. Values sourced from latest Universalis
# snapshot; re-derive currentAveragePrice* / LevePrice* / LeveProfit*.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 141
$ws.Range("I12").Value = 141
$ws.Range("K12").Value = 141
$ws.Range("M12").Value = 29
# Row 28
$ws.Range("H28").Value = 453.13635
$ws.Range("I28").Value = 403.6316
$ws.Range("J28").Value = 766.6667
$ws.Range("K28").Value = 403.6316
$ws.Range("L28").Value = 766.6667
$ws.Range("M28").Value = 81.36840000000001
$ws.Range("N28").Value = -1736.6667
# Row 33
$ws.Range("H33").Value = 941.95654
$ws.Range("I33").Value = 1096.2858
$ws.Range("J33").Value = 701.8889
$ws.Range("K33").Value = 1096.2858
$ws.Range("L33").Value = 701.8889
$ws.Range("M33").Value = -867.2858000000001
$ws.Range("N33").Value = -1159.8889
# Row 98
$ws.Range("H98").Value = 1528.6
$ws.Range("I98").Value = 1622.5714
$ws.Range("J98").Value = 1409
$ws.Range("K98").Value = 1622.5714
$ws.Range("L98").Value = 1409
$ws.Range("M98").Value = -124.5714
$ws.Range("N98").Value = -4405
# Row 122
$ws.Range("H122").Value = 1528.6
$ws.Range("I122").Value = 1622.5714
$ws.Range("J122").Value = 1409
$ws.Range("K122").Value = 4867.7142
$ws.Range("L122").Value = 4227
$ws.Range("M122").Value = -2417.7142
$ws.Range("N122").Value = -9127
# Row 138
$ws.Range("H138").Value = 1936.4637
$ws.Range("I138").Value = 968
$ws.Range("J138").Value = 3442.963
$ws.Range("K138").Value = 2904
$ws.Range("L138").Value = 10328.889
$ws.Range("M138").Value = 2236
$ws.Range("N138").Value = -20608.889

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17895.111
$ws.Range("I32").Value = 18419.13
$ws.Range("K32").Value = 18419.13
$ws.Range("M32").Value = -18132.13
# Row 88
$ws.Range("H88").Value = 1959.2222
$ws.Range("I88").Value = 1430.6
$ws.Range("J88").Value = 2620
$ws.Range("K88").Value = 1430.6
$ws.Range("L88").Value = 2620
$ws.Range("M88").Value = -1024.6
$ws.Range("N88").Value = -3432
# Row 91
$ws.Range("H91").Value = 1959.2222
$ws.Range("I91").Value = 1430.6
$ws.Range("J91").Value = 2620
$ws.Range("K91").Value = 1430.6
$ws.Range("L91").Value = 2620
$ws.Range("M91").Value = -26.59999999999991
$ws.Range("N91").Value = -5428
# Row 114
$ws.Range("H114").Value = 33050
$ws.Range("J114").Value = 33050
$ws.Range("L114").Value = 33050
$ws.Range("N114").Value = -41728
# Row 121
$ws.Range("H121").Value = 19127.5
$ws.Range("J121").Value = 19127.5
$ws.Range("L121").Value = 19127.5
$ws.Range("N121").Value = -22621.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 44
$ws.Range("H44").Value = 19500
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 86
$ws.Range("H86").Value = 1356.2858
$ws.Range("I86").Value = 1298.8
$ws.Range("K86").Value = 1298.8
$ws.Range("M86").Value = -175.8
# Row 89
$ws.Range("H89").Value = 1356.2858
$ws.Range("I89").Value = 1298.8
$ws.Range("K89").Value = 6122.5
$ws.Range("M89").Value = -878

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 30218.25
$ws.Range("I62").Value = 35248.668
$ws.Range("J62").Value = 27200
$ws.Range("K62").Value = 35248.668
$ws.Range("L62").Value = 27200
$ws.Range("M62").Value = -34624.668
$ws.Range("N62").Value = -28448
# Row 65
$ws.Range("H65").Value = 30218.25
$ws.Range("I65").Value = 35248.668
$ws.Range("J65").Value = 27200
$ws.Range("K65").Value = 176243.34
$ws.Range("L65").Value = 136000
$ws.Range("M65").Value = -173123.34
$ws.Range("N65").Value = -142240
# Row 99
$ws.Range("H99").Value = 2241.348
$ws.Range("I99").Value = 1901.7142
$ws.Range("J99").Value = 2769.6667
$ws.Range("K99").Value = 1901.7142
$ws.Range("L99").Value = 2769.6667
$ws.Range("M99").Value = -403.7141999999999
$ws.Range("N99").Value = -5765.6667
# Row 126
$ws.Range("H126").Value = 2241.348
$ws.Range("I126").Value = 1901.7142
$ws.Range("J126").Value = 2769.6667
$ws.Range("K126").Value = 5705.142599999999
$ws.Range("L126").Value = 8309.000100000001
$ws.Range("M126").Value = -3235.142599999999
$ws.Range("N126").Value = -13249.0001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 3166.6667
$ws.Range("I36").Value = 750
$ws.Range("J36").Value = 8000
$ws.Range("K36").Value = 2250
$ws.Range("L36").Value = 24000
$ws.Range("M36").Value = -2081
$ws.Range("N36").Value = -24338
# Row 132
$ws.Range("H132").Value = 2706.3333
$ws.Range("I132").Value = 1057
$ws.Range("J132").Value = 6005
$ws.Range("K132").Value = 9513
$ws.Range("L132").Value = 54045
$ws.Range("M132").Value = -6983
$ws.Range("N132").Value = -59105

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 46250
$ws.Range("J51").Value = 46250
$ws.Range("L51").Value = 46250
$ws.Range("N51").Value = -47268
# Row 113
$ws.Range("H113").Value = 2870
$ws.Range("I113").Value = 2475
$ws.Range("J113").Value = 3186
$ws.Range("K113").Value = 2475
$ws.Range("L113").Value = 3186
$ws.Range("M113").Value = -305
$ws.Range("N113").Value = -7526
# Row 122
$ws.Range("H122").Value = 2172.697
$ws.Range("I122").Value = 2125.8333
$ws.Range("J122").Value = 2297.6667
$ws.Range("K122").Value = 6377.499899999999
$ws.Range("L122").Value = 6893.000100000001
$ws.Range("M122").Value = -3927.499899999999
$ws.Range("N122").Value = -11793.0001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1004.82355
$ws.Range("I22").Value = 732.5
$ws.Range("J22").Value = 1246.8889
$ws.Range("K22").Value = 732.5
$ws.Range("L22").Value = 1246.8889
$ws.Range("M22").Value = -437.5
$ws.Range("N22").Value = -1836.8889
# Row 27
$ws.Range("H27").Value = 1004.82355
$ws.Range("I27").Value = 732.5
$ws.Range("J27").Value = 1246.8889
$ws.Range("K27").Value = 732.5
$ws.Range("L27").Value = 1246.8889
$ws.Range("M27").Value = -625.5
$ws.Range("N27").Value = -1460.8889
# Row 68
$ws.Range("H68").Value = 3788.3333
$ws.Range("I68").Value = 3626.4
$ws.Range("J68").Value = 3990.75
$ws.Range("K68").Value = 3626.4
$ws.Range("L68").Value = 3990.75
$ws.Range("M68").Value = -2877.4
$ws.Range("N68").Value = -5488.75
# Row 71
$ws.Range("H71").Value = 3788.3333
$ws.Range("I71").Value = 3626.4
$ws.Range("J71").Value = 3990.75
$ws.Range("K71").Value = 18132
$ws.Range("L71").Value = 19953.75
$ws.Range("M71").Value = -14388
$ws.Range("N71").Value = -27441.75
# Row 103
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
# Row 132
$ws.Range("H132").Value = 21279056
$ws.Range("I132").Value = 30305344
$ws.Range("J132").Value = 2807.3572
$ws.Range("K132").Value = 90916032
$ws.Range("L132").Value = 8422.071599999999
$ws.Range("M132").Value = -90913502
$ws.Range("N132").Value = -13482.0716

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 530.7692
$ws.Range("I14").Value = 508.33334
$ws.Range("K14").Value = 508.33334
$ws.Range("M14").Value = -340.33334
